$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

function Find-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        if ($shapes.Item($i).Id -eq $id) {
            return $shapes.Item($i)
        }
    }
    return $null
}

# --- Group the connector + label pairs (mirrors manual grouping in PowerPoint) ---
$conn10 = Find-ShapeById $s.Shapes 11
$lbl247 = Find-ShapeById $s.Shapes 23
$grp1 = $s.Shapes.Range(@($conn10.Name, $lbl247.Name)).Group()

$conn11 = Find-ShapeById $s.Shapes 12
$lbl895 = Find-ShapeById $s.Shapes 24
$grp2 = $s.Shapes.Range(@($conn11.Name, $lbl895.Name)).Group()

$conn12 = Find-ShapeById $s.Shapes 13
$lbl73 = Find-ShapeById $s.Shapes 25
$grp3 = $s.Shapes.Range(@($conn12.Name, $lbl73.Name)).Group()

# --- Animate the three new groups with a Wipe (from top) entrance, first on click, rest with previous ---
$seq = $s.TimeLine.MainSequence

$e1 = $seq.AddEffect($grp3, 22, 0, 1)
$e1.EffectParameters.Direction = 3
$e1.Timing.Duration = 1.3
$e1.Timing.RepeatCount = 4

$e2 = $seq.AddEffect($grp2, 22, 0, 2)
$e2.EffectParameters.Direction = 3
$e2.Timing.Duration = 1.25
$e2.Timing.RepeatCount = 4

$e3 = $seq.AddEffect($grp1, 22, 0, 2)
$e3.EffectParameters.Direction = 3
$e3.Timing.Duration = 1.25
$e3.Timing.RepeatCount = 4
